$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "331.36"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.18%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "11"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "45.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.29%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "11"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.550"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.00%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "11"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08345"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.88%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "11"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.089"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.95%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "11"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9829"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.09%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "11"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.538"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.85%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "11"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1203"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5.02%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "11"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1923"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.83%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "11"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "10.31"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.19%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "11"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09892"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.49%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "11"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04670"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.14%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "11"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1059"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.44%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "11"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001288"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.24%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "11"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005943"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.97%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "11"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.392"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.79%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "11"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.456"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "11"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3349"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-3.17%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "11"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1355"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.57%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "11"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2567"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.52%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "11"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04170"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.39%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "11"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.58%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "11"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004581"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "5.06%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "11"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001302"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.57%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "11"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003749"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.19%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "11"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "11"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "11"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "11"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "11"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "11"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "11"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "11"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "11"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "11"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "11"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "11"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02689"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "3.71%"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "11"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.39%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "11"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007908"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.62%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "11"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1430"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.89%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "11"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007518"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "11"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002121"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "5.31%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "11"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008476"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.66%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "11"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3398"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "11"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007102"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.68%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "11"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.30%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "11"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.37%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "11"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003531"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-4.11%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "11"
$ws.Range("B50").Value = "CoinbaseStockToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.003538"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.30%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "11"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00002105"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.30%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "11"
